$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3: paper found via CrossRef lookup (previously 'Unknown Title' / not found) ---
$ws.Range("C3").Value = '"Characteristics and Outcomes of 21 Critically Ill Patients With COVID-19 in Washington State"'
$ws.Range("E3").Value = '[Matt%Arentz%xref no email%0,  Eric%Yim%xref no email%0,  Lindy%Klaff%xref no email%0,  Sharukh%Lokhandwala%xref no email%0,  Francis X.%Riedo%xref no email%0,  Maria%Chong%xref no email%0,  Melissa%Lee%xref no email%0]'
$ws.Range("F3").Value = '10.1001/jama.2020.4326'
$ws.Range("G3").Value = 'CROSSREF'

# H3 holds a literal text date (not a real Excel date serial). Assigning a
# "yyyy-mm-dd" string straight to .Value gets auto-parsed into a date serial,
# so stage it as a text formula in a scratch cell, then paste-special just the
# resulting value into H3. This keeps it a plain text/shared-string cell and
# avoids creating any new (unused) number-format style.
$scratch = $ws.Range("Z1")
$scratch.Formula = '="2023-05-25"'
$scratch.Copy()
$ws.Range("H3").PasteSpecial(-4163)
$scratch.Clear()

# --- Rows 4 & 5: author-name lists re-synced (re-fetched with extra padding) ---
$ws.Range("E4").Value = '[Ruchong%Chen%NULL%0,         Wenhua%Liang%NULL%0,         Mei%Jiang%NULL%0,         Weijie%Guan%NULL%0,         Chen%Zhan%NULL%0,         Tao%Wang%NULL%0,         Chunli%Tang%NULL%0,         Ling%Sang%NULL%0,         Jiaxing%Liu%NULL%0,         Zhengyi%Ni%NULL%0,         Yu%Hu%NULL%0,         Lei%Liu%NULL%0,         Hong%Shan%NULL%0,         Chunliang%Lei%NULL%0,         Yixiang%Peng%NULL%0,         Li%Wei%NULL%0,         Yong%Liu%NULL%0,         Yahua%Hu%NULL%0,         Peng%Peng%NULL%0,         Jianming%Wang%NULL%0,         Jiyang%Liu%NULL%0,         Zhong%Chen%NULL%0,         Gang%Li%NULL%0,         Zhijian%Zheng%NULL%0,         Shaoqin%Qiu%NULL%0,         Jie%Luo%NULL%0,         Changjiang%Ye%NULL%0,         Shaoyong%Zhu%NULL%0,         Xiaoqing%Liu%NULL%0,         Linling%Cheng%NULL%0,         Feng%Ye%NULL%0,         Jinping%Zheng%NULL%0,         Nuofu%Zhang%NULL%0,         Yimin%Li%NULL%0,         Jianxing%He%NULL%0,         Shiyue%Li%lishiyue@188.com%0,         Nanshan%Zhong%NULL%0,         NULL%NULL%NULL%0]'
$ws.Range("E5").Value = '[Wei-jie%Guan%NULL%0,         Zheng-yi%Ni%NULL%6,         Zheng-yi%Ni%NULL%0,         Yu%Hu%NULL%0,         Wen-hua%Liang%NULL%3,         Chun-quan%Ou%NULL%3,         Jian-xing%He%NULL%3,         Lei%Liu%NULL%0,         Hong%Shan%NULL%0,         Chun-liang%Lei%NULL%3,         David S.C.%Hui%NULL%4,         Bin%Du%NULL%3,         Lan-juan%Li%NULL%3,         Guang%Zeng%NULL%3,         Kwok-Yung%Yuen%NULL%0,         Ru-chong%Chen%NULL%3,         Chun-li%Tang%NULL%3,         Tao%Wang%NULL%0,         Ping-yan%Chen%NULL%3,         Jie%Xiang%NULL%0,         Shi-yue%Li%NULL%3,         Jin-lin%Wang%NULL%3,         Zi-jing%Liang%NULL%3,         Yi-xiang%Peng%NULL%3,         Li%Wei%NULL%0,         Yong%Liu%NULL%0,         Ya-hua%Hu%NULL%3,         Peng%Peng%NULL%0,         Jian-ming%Wang%NULL%4,         Ji-yang%Liu%NULL%3,         Zhong%Chen%NULL%0,         Gang%Li%NULL%0,         Zhi-jian%Zheng%NULL%3,         Shao-qin%Qiu%NULL%3,         Jie%Luo%NULL%0,         Chang-jiang%Ye%NULL%3,         Shao-yong%Zhu%NULL%3,         Nan-shan%Zhong%NULL%3]'

